$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Unmerge the blocks whose row-span changes (Event, Language, Organization, Species)
$ws.Range("A16:A20").UnMerge()
$ws.Range("A21:A25").UnMerge()
$ws.Range("A26:A30").UnMerge()
$ws.Range("A37:A41").UnMerge()

# Remove the now-extra last row (table shrinks from 41 to 40 data+header rows)
$ws.Rows.Item(41).Delete()

# Rewrite every data cell (rows 2-40) to the updated / regenerated template content
# Row 2: Country
$ws.Cells.Item(2, 1).Value = "Country"
$ws.Cells.Item(2, 2).Value = "What is the ISO alpha-2 code for {country}?"
$ws.Cells.Item(2, 3).Value = 36
$ws.Cells.Item(2, 4).Value = 1
$ws.Cells.Item(2, 5).Value = 1.027777777777778
$ws.Cells.Item(2, 6).Value = "What is the ISO alpha-2 code for Nigeria?"
$ws.Cells.Item(2, 7).Value = "NG"
$ws.Cells.Item(2, 8).Value = "NG"

# Row 3: 
$ws.Cells.Item(3, 1).Value = ""
$ws.Cells.Item(3, 2).Value = "What is the calling code for {country}?"
$ws.Cells.Item(3, 3).Value = 36
$ws.Cells.Item(3, 4).Value = 1
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = "What is the calling code for Armenia?"
$ws.Cells.Item(3, 7).Value = "+374"
$ws.Cells.Item(3, 8).Value = "+374"

# Row 4: 
$ws.Cells.Item(4, 1).Value = ""
$ws.Cells.Item(4, 2).Value = "What is the capital of {country}?"
$ws.Cells.Item(4, 3).Value = 36
$ws.Cells.Item(4, 4).Value = 1
$ws.Cells.Item(4, 5).Value = 2.277777777777778
$ws.Cells.Item(4, 6).Value = "What is the capital of Turkey?"
$ws.Cells.Item(4, 7).Value = "Ankara"
$ws.Cells.Item(4, 8).Value = "Ankara"

# Row 5: 
$ws.Cells.Item(5, 1).Value = ""
$ws.Cells.Item(5, 2).Value = "What is the currency of {country}?"
$ws.Cells.Item(5, 3).Value = 36
$ws.Cells.Item(5, 4).Value = 0.9833333333333333
$ws.Cells.Item(5, 5).Value = 3.611111111111111
$ws.Cells.Item(5, 6).Value = "What is the currency of Vietnam?"
$ws.Cells.Item(5, 7).Value = "Vietnamese dong"
$ws.Cells.Item(5, 8).Value = "Dong"

# Row 6: 
$ws.Cells.Item(6, 1).Value = ""
$ws.Cells.Item(6, 2).Value = "What is the top-level internet domain for {country}?"
$ws.Cells.Item(6, 3).Value = 36
$ws.Cells.Item(6, 4).Value = 1
$ws.Cells.Item(6, 5).Value = 1.138888888888889
$ws.Cells.Item(6, 6).Value = "What is the top-level internet domain for Spain?"
$ws.Cells.Item(6, 7).Value = ".es"
$ws.Cells.Item(6, 8).Value = ".es"

# Row 7: 
$ws.Cells.Item(7, 1).Value = ""
$ws.Cells.Item(7, 2).Value = "What language in {country} has the most speakers?"
$ws.Cells.Item(7, 3).Value = 36
$ws.Cells.Item(7, 4).Value = 0.9916666666666667
$ws.Cells.Item(7, 5).Value = 2.25
$ws.Cells.Item(7, 6).Value = "What language in Kenya has the most speakers?"
$ws.Cells.Item(7, 7).Value = "Swahili"
$ws.Cells.Item(7, 8).Value = "Swahili"

# Row 8: 
$ws.Cells.Item(8, 1).Value = ""
$ws.Cells.Item(8, 2).Value = "Which ethnic group is the largest in {country}?"
$ws.Cells.Item(8, 3).Value = 36
$ws.Cells.Item(8, 4).Value = 0.9750000000000001
$ws.Cells.Item(8, 5).Value = 2.694444444444445
$ws.Cells.Item(8, 6).Value = "Which ethnic group is the largest in Thailand?"
$ws.Cells.Item(8, 7).Value = "Thai"
$ws.Cells.Item(8, 8).Value = "Thai"

# Row 9: 
$ws.Cells.Item(9, 1).Value = ""
$ws.Cells.Item(9, 2).Value = "Which religion has the most followers in {country}?"
$ws.Cells.Item(9, 3).Value = 36
$ws.Cells.Item(9, 4).Value = 0.9750000000000001
$ws.Cells.Item(9, 5).Value = 2.388888888888889
$ws.Cells.Item(9, 6).Value = "Which religion has the most followers in Norway?"
$ws.Cells.Item(9, 7).Value = "Christianity"
$ws.Cells.Item(9, 8).Value = "Christianity"

# Row 10: Creative Work
$ws.Cells.Item(10, 1).Value = "Creative Work"
$ws.Cells.Item(10, 2).Value = "In which country was {creative_work} first released or published?"
$ws.Cells.Item(10, 3).Value = 21
$ws.Cells.Item(10, 4).Value = 1
$ws.Cells.Item(10, 5).Value = 1.619047619047619
$ws.Cells.Item(10, 6).Value = "In which country was Amélie first released or published?"
$ws.Cells.Item(10, 7).Value = "France"
$ws.Cells.Item(10, 8).Value = "France"

# Row 11: 
$ws.Cells.Item(11, 1).Value = ""
$ws.Cells.Item(11, 2).Value = "What is the genre or style of {creative_work}?"
$ws.Cells.Item(11, 3).Value = 21
$ws.Cells.Item(11, 4).Value = 0.7904761904761903
$ws.Cells.Item(11, 5).Value = 4.761904761904762
$ws.Cells.Item(11, 6).Value = "What is the genre or style of The Brothers Karamazov?"
$ws.Cells.Item(11, 7).Value = "Philosophical, psychological fiction"
$ws.Cells.Item(11, 8).Value = "Psychological drama"

# Row 12: 
$ws.Cells.Item(12, 1).Value = ""
$ws.Cells.Item(12, 2).Value = "What is the original language of {creative_work}?"
$ws.Cells.Item(12, 3).Value = 21
$ws.Cells.Item(12, 4).Value = 1
$ws.Cells.Item(12, 5).Value = 1.047619047619048
$ws.Cells.Item(12, 6).Value = "What is the original language of The Grapes of Wrath?"
$ws.Cells.Item(12, 7).Value = "English"
$ws.Cells.Item(12, 8).Value = "English"

# Row 13: 
$ws.Cells.Item(13, 1).Value = ""
$ws.Cells.Item(13, 2).Value = "When was {creative_work} released or published?"
$ws.Cells.Item(13, 3).Value = 21
$ws.Cells.Item(13, 4).Value = 0.9761904761904762
$ws.Cells.Item(13, 5).Value = 2.476190476190476
$ws.Cells.Item(13, 6).Value = "When was A Separation released or published?"
$ws.Cells.Item(13, 7).Value = "2011"
$ws.Cells.Item(13, 8).Value = "2011"

# Row 14: 
$ws.Cells.Item(14, 1).Value = ""
$ws.Cells.Item(14, 2).Value = "Where was {creative_work} produced or created?"
$ws.Cells.Item(14, 3).Value = 21
$ws.Cells.Item(14, 4).Value = 0.9571428571428571
$ws.Cells.Item(14, 5).Value = 1.952380952380952
$ws.Cells.Item(14, 6).Value = "Where was War and Peace produced or created?"
$ws.Cells.Item(14, 7).Value = "Russia"
$ws.Cells.Item(14, 8).Value = "Russia"

# Row 15: 
$ws.Cells.Item(15, 1).Value = ""
$ws.Cells.Item(15, 2).Value = "Who is the creator of {creative_work}?"
$ws.Cells.Item(15, 3).Value = 21
$ws.Cells.Item(15, 4).Value = 0.9857142857142857
$ws.Cells.Item(15, 5).Value = 3.952380952380953
$ws.Cells.Item(15, 6).Value = "Who is the creator of Run Lola Run?"
$ws.Cells.Item(15, 7).Value = "Tom Tykwer"
$ws.Cells.Item(15, 8).Value = "Tom Tykwer"

# Row 16: Event
$ws.Cells.Item(16, 1).Value = "Event"
$ws.Cells.Item(16, 2).Value = "In which country did {event} happen?"
$ws.Cells.Item(16, 3).Value = 35
$ws.Cells.Item(16, 4).Value = 0.9857142857142858
$ws.Cells.Item(16, 5).Value = 1.485714285714286
$ws.Cells.Item(16, 6).Value = "In which country did French Revolution happen?"
$ws.Cells.Item(16, 7).Value = "France"
$ws.Cells.Item(16, 8).Value = "France"

# Row 17: 
$ws.Cells.Item(17, 1).Value = ""
$ws.Cells.Item(17, 2).Value = "What year did {event} end?"
$ws.Cells.Item(17, 3).Value = 35
$ws.Cells.Item(17, 4).Value = 0.9857142857142858
$ws.Cells.Item(17, 5).Value = 2.114285714285714
$ws.Cells.Item(17, 6).Value = "What year did American Civil War end?"
$ws.Cells.Item(17, 7).Value = "1865"
$ws.Cells.Item(17, 8).Value = "1865"

# Row 18: 
$ws.Cells.Item(18, 1).Value = ""
$ws.Cells.Item(18, 2).Value = "When did {event} take place?"
$ws.Cells.Item(18, 3).Value = 35
$ws.Cells.Item(18, 4).Value = 0.9342857142857144
$ws.Cells.Item(18, 5).Value = 4.914285714285715
$ws.Cells.Item(18, 6).Value = "When did The Surrender of Japan in WWII take place?"
$ws.Cells.Item(18, 7).Value = "September 2, 1945"
$ws.Cells.Item(18, 8).Value = "September 2, 1945"

# Row 19: 
$ws.Cells.Item(19, 1).Value = ""
$ws.Cells.Item(19, 2).Value = "Who was the most important leader or figure involved in {event}?"
$ws.Cells.Item(19, 3).Value = 35
$ws.Cells.Item(19, 4).Value = 0.9914285714285715
$ws.Cells.Item(19, 5).Value = 4.371428571428571
$ws.Cells.Item(19, 6).Value = "Who was the most important leader or figure involved in French Revolution?"
$ws.Cells.Item(19, 7).Value = "Maximilien Robespierre"
$ws.Cells.Item(19, 8).Value = "Robespierre"

# Row 20: Language
$ws.Cells.Item(20, 1).Value = "Language"
$ws.Cells.Item(20, 2).Value = "What is the ISO 639‑1 code for {language}?"
$ws.Cells.Item(20, 3).Value = 21
$ws.Cells.Item(20, 4).Value = 1
$ws.Cells.Item(20, 5).Value = 1
$ws.Cells.Item(20, 6).Value = "What is the ISO 639‑1 code for Kazakh?"
$ws.Cells.Item(20, 7).Value = "kk"
$ws.Cells.Item(20, 8).Value = "kk"

# Row 21: 
$ws.Cells.Item(21, 1).Value = ""
$ws.Cells.Item(21, 2).Value = "What is the name of the alphabet or script of {language}?"
$ws.Cells.Item(21, 3).Value = 21
$ws.Cells.Item(21, 4).Value = 1
$ws.Cells.Item(21, 5).Value = 2.619047619047619
$ws.Cells.Item(21, 6).Value = "What is the name of the alphabet or script of Greek?"
$ws.Cells.Item(21, 7).Value = "Greek alphabet"
$ws.Cells.Item(21, 8).Value = "Greek alphabet"

# Row 22: 
$ws.Cells.Item(22, 1).Value = ""
$ws.Cells.Item(22, 2).Value = "What is the primary word order in {language}?"
$ws.Cells.Item(22, 3).Value = 21
$ws.Cells.Item(22, 4).Value = 1
$ws.Cells.Item(22, 5).Value = 8.619047619047619
$ws.Cells.Item(22, 6).Value = "What is the primary word order in Haitian Creole?"
$ws.Cells.Item(22, 7).Value = "SVO (Subject-Verb-Object)"
$ws.Cells.Item(22, 8).Value = "SVO (Subject-Verb-Object)"

# Row 23: 
$ws.Cells.Item(23, 1).Value = ""
$ws.Cells.Item(23, 2).Value = "What region is {language} native to?"
$ws.Cells.Item(23, 3).Value = 21
$ws.Cells.Item(23, 4).Value = 0.8619047619047618
$ws.Cells.Item(23, 5).Value = 3.666666666666667
$ws.Cells.Item(23, 6).Value = "What region is Kazakh native to?"
$ws.Cells.Item(23, 7).Value = "Central Asia"
$ws.Cells.Item(23, 8).Value = "Central Asia"

# Row 24: 
$ws.Cells.Item(24, 1).Value = ""
$ws.Cells.Item(24, 2).Value = "What writing system is used by {language}?"
$ws.Cells.Item(24, 3).Value = 21
$ws.Cells.Item(24, 4).Value = 1
$ws.Cells.Item(24, 5).Value = 2.761904761904762
$ws.Cells.Item(24, 6).Value = "What writing system is used by Haitian Creole?"
$ws.Cells.Item(24, 7).Value = "Latin alphabet"
$ws.Cells.Item(24, 8).Value = "Latin alphabet"

# Row 25: Organization
$ws.Cells.Item(25, 1).Value = "Organization"
$ws.Cells.Item(25, 2).Value = "In what year was {organization} established?"
$ws.Cells.Item(25, 3).Value = 22
$ws.Cells.Item(25, 4).Value = 1
$ws.Cells.Item(25, 5).Value = 2
$ws.Cells.Item(25, 6).Value = "In what year was Alibaba established?"
$ws.Cells.Item(25, 7).Value = "1999"
$ws.Cells.Item(25, 8).Value = "1999"

# Row 26: 
$ws.Cells.Item(26, 1).Value = ""
$ws.Cells.Item(26, 2).Value = "What is the primary field or industry of {organization}?"
$ws.Cells.Item(26, 3).Value = 22
$ws.Cells.Item(26, 4).Value = 0.8909090909090908
$ws.Cells.Item(26, 5).Value = 2.772727272727273
$ws.Cells.Item(26, 6).Value = "What is the primary field or industry of Airbnb?"
$ws.Cells.Item(26, 7).Value = "Hospitality"
$ws.Cells.Item(26, 8).Value = "Short-term rentals"

# Row 27: 
$ws.Cells.Item(27, 1).Value = ""
$ws.Cells.Item(27, 2).Value = "What primary service or product does {organization} provide?"
$ws.Cells.Item(27, 3).Value = 22
$ws.Cells.Item(27, 4).Value = 0.8681818181818181
$ws.Cells.Item(27, 5).Value = 3.727272727272727
$ws.Cells.Item(27, 6).Value = "What primary service or product does Red Cross provide?"
$ws.Cells.Item(27, 7).Value = "Emergency humanitarian aid"
$ws.Cells.Item(27, 8).Value = "Disaster relief and humanitarian aid"

# Row 28: 
$ws.Cells.Item(28, 1).Value = ""
$ws.Cells.Item(28, 2).Value = "Where is the headquarters of {organization} located?"
$ws.Cells.Item(28, 3).Value = 22
$ws.Cells.Item(28, 4).Value = 1
$ws.Cells.Item(28, 5).Value = 5.136363636363637
$ws.Cells.Item(28, 6).Value = "Where is the headquarters of Alibaba located?"
$ws.Cells.Item(28, 7).Value = "Hangzhou, China"
$ws.Cells.Item(28, 8).Value = "Hangzhou, China"

# Row 29: 
$ws.Cells.Item(29, 1).Value = ""
$ws.Cells.Item(29, 2).Value = "Where was {organization} established?"
$ws.Cells.Item(29, 3).Value = 22
$ws.Cells.Item(29, 4).Value = 0.959090909090909
$ws.Cells.Item(29, 5).Value = 4.272727272727272
$ws.Cells.Item(29, 6).Value = "Where was Johnson & Johnson established?"
$ws.Cells.Item(29, 7).Value = "New Brunswick, New Jersey, USA"
$ws.Cells.Item(29, 8).Value = "New Brunswick, New Jersey"

# Row 30: 
$ws.Cells.Item(30, 1).Value = ""
$ws.Cells.Item(30, 2).Value = "Who established {organization}?"
$ws.Cells.Item(30, 3).Value = 22
$ws.Cells.Item(30, 4).Value = 0.8136363636363636
$ws.Cells.Item(30, 5).Value = 5.363636363636363
$ws.Cells.Item(30, 6).Value = "Who established Coca-Cola?"
$ws.Cells.Item(30, 7).Value = "John Stith Pemberton"
$ws.Cells.Item(30, 8).Value = "John Pemberton"

# Row 31: Person
$ws.Cells.Item(31, 1).Value = "Person"
$ws.Cells.Item(31, 2).Value = "What language was primarily spoken by {person}?"
$ws.Cells.Item(31, 3).Value = 26
$ws.Cells.Item(31, 4).Value = 0.9769230769230769
$ws.Cells.Item(31, 5).Value = 1.423076923076923
$ws.Cells.Item(31, 6).Value = "What language was primarily spoken by William Shakespeare?"
$ws.Cells.Item(31, 7).Value = "Early Modern English"
$ws.Cells.Item(31, 8).Value = "English"

# Row 32: 
$ws.Cells.Item(32, 1).Value = ""
$ws.Cells.Item(32, 2).Value = "What occupation is {person} most well-known for?"
$ws.Cells.Item(32, 3).Value = 26
$ws.Cells.Item(32, 4).Value = 0.9653846153846153
$ws.Cells.Item(32, 5).Value = 3.153846153846154
$ws.Cells.Item(32, 6).Value = "What occupation is Alexander the Great most well-known for?"
$ws.Cells.Item(32, 7).Value = "Military leader and king"
$ws.Cells.Item(32, 8).Value = "Conqueror and military leader"

# Row 33: 
$ws.Cells.Item(33, 1).Value = ""
$ws.Cells.Item(33, 2).Value = "What year did {person} pass away?"
$ws.Cells.Item(33, 3).Value = 26
$ws.Cells.Item(33, 4).Value = 1
$ws.Cells.Item(33, 5).Value = 2
$ws.Cells.Item(33, 6).Value = "What year did Cleopatra VII pass away?"
$ws.Cells.Item(33, 7).Value = "30 BC"
$ws.Cells.Item(33, 8).Value = "30 BC"

# Row 34: 
$ws.Cells.Item(34, 1).Value = ""
$ws.Cells.Item(34, 2).Value = "What year was {person} born?"
$ws.Cells.Item(34, 3).Value = 26
$ws.Cells.Item(34, 4).Value = 1
$ws.Cells.Item(34, 5).Value = 2.115384615384615
$ws.Cells.Item(34, 6).Value = "What year was Alexander the Great born?"
$ws.Cells.Item(34, 7).Value = "356 BC"
$ws.Cells.Item(34, 8).Value = "356 BC"

# Row 35: 
$ws.Cells.Item(35, 1).Value = ""
$ws.Cells.Item(35, 2).Value = "Where did {person} die?"
$ws.Cells.Item(35, 3).Value = 26
$ws.Cells.Item(35, 4).Value = 0.9307692307692309
$ws.Cells.Item(35, 5).Value = 5.076923076923077
$ws.Cells.Item(35, 6).Value = "Where did Vladimir Lenin die?"
$ws.Cells.Item(35, 7).Value = "Gorki, Russia"
$ws.Cells.Item(35, 8).Value = "Dacha at Gorky Park, Moscow, Russia"

# Row 36: 
$ws.Cells.Item(36, 1).Value = ""
$ws.Cells.Item(36, 2).Value = "Where was the birthplace of {person}?"
$ws.Cells.Item(36, 3).Value = 26
$ws.Cells.Item(36, 4).Value = 0.9538461538461539
$ws.Cells.Item(36, 5).Value = 5.153846153846154
$ws.Cells.Item(36, 6).Value = "Where was the birthplace of Vincent van Gogh?"
$ws.Cells.Item(36, 7).Value = "Zundert, Netherlands"
$ws.Cells.Item(36, 8).Value = "Groot-Zundert, Netherlands"

# Row 37: Species
$ws.Cells.Item(37, 1).Value = "Species"
$ws.Cells.Item(37, 2).Value = "What is the diet of {species}?"
$ws.Cells.Item(37, 3).Value = 25
$ws.Cells.Item(37, 4).Value = 0.768
$ws.Cells.Item(37, 5).Value = 10.32
$ws.Cells.Item(37, 6).Value = "What is the diet of wolverine?"
$ws.Cells.Item(37, 7).Value = "Carnivorous; eats small animals, birds, carrion, and berries"
$ws.Cells.Item(37, 8).Value = "Carnivorous, eating small mammals, birds, and reptiles"

# Row 38: 
$ws.Cells.Item(38, 1).Value = ""
$ws.Cells.Item(38, 2).Value = "What is the social structure of {species}?"
$ws.Cells.Item(38, 3).Value = 25
$ws.Cells.Item(38, 4).Value = 0.772
$ws.Cells.Item(38, 5).Value = 6.12
$ws.Cells.Item(38, 6).Value = "What is the social structure of tiger?"
$ws.Cells.Item(38, 7).Value = "Solitary, territorial"
$ws.Cells.Item(38, 8).Value = "Solitary or small groups"

# Row 39: 
$ws.Cells.Item(39, 1).Value = ""
$ws.Cells.Item(39, 2).Value = "What type of organism is this {species}?"
$ws.Cells.Item(39, 3).Value = 25
$ws.Cells.Item(39, 4).Value = 0.968
$ws.Cells.Item(39, 5).Value = 3.08
$ws.Cells.Item(39, 6).Value = "What type of organism is this panda?"
$ws.Cells.Item(39, 7).Value = "Mammal"
$ws.Cells.Item(39, 8).Value = "Mammal"

# Row 40: 
$ws.Cells.Item(40, 1).Value = ""
$ws.Cells.Item(40, 2).Value = "Where is {species} primarily native to?"
$ws.Cells.Item(40, 3).Value = 25
$ws.Cells.Item(40, 4).Value = 0.848
$ws.Cells.Item(40, 5).Value = 3.68
$ws.Cells.Item(40, 6).Value = "Where is red-shouldered hawk primarily native to?"
$ws.Cells.Item(40, 7).Value = "North America"
$ws.Cells.Item(40, 8).Value = "North America"

# Re-merge column A for each entity-type block with its new row span
$ws.Range("A16:A19").Merge()
$ws.Range("A20:A24").Merge()
$ws.Range("A25:A30").Merge()
$ws.Range("A37:A40").Merge()

$ws.Range("A1").Select()
